$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add 2 hours worked on Friday (3/4) for the "Task zzzzzz" row
$ws.Range("G13").Value = 2

# Update the active selection to match the edited cell (H13), matching author's final cursor position
$ws.Range("H13").Select()
